# 036 : added exctsccr driver
# Moves the 4 "Exciting Soccer" (exctsccr.c) rows out of the "ALL" master
# list and into the per-status sheets: the 3 playable sets go to
# "Playable (untested)", and the non-working "exctscc2" clone goes to
# "GAME_NOT_WORKING FLAG".

$wb  = $excel.ActiveWorkbook
$all = $wb.Worksheets.Item("ALL")
$playable = $wb.Worksheets.Item("Playable (untested)")
$notWorking = $wb.Worksheets.Item("GAME_NOT_WORKING FLAG")

# --- append the 3 playable rows to "Playable (untested)" (rows 295-297) ---
$playable.Range("A295").Value = 295
$playable.Range("B295").Value = "exctsccr"
$playable.Range("C295").Value = "exctsccr.c"
$playable.Range("D295").Value = "Z80"
$playable.Range("E295").Value = "Z80"
$playable.Range("H295").Value = "4xAY-8910"
$playable.Range("I295").Value = "2xDAC"
$playable.Range("M295").Value = "Exciting Soccer"

$playable.Range("A296").Value = 296
$playable.Range("B296").Value = "exctscca"
$playable.Range("C296").Value = "exctsccr.c"
$playable.Range("D296").Value = "Z80"
$playable.Range("E296").Value = "Z80"
$playable.Range("H296").Value = "4xAY-8910"
$playable.Range("I296").Value = "2xDAC"
$playable.Range("M296").Value = "Exciting Soccer (alternate music)"

$playable.Range("A297").Value = 297
$playable.Range("B297").Value = "exctsccb"
$playable.Range("C297").Value = "exctsccr.c"
$playable.Range("D297").Value = "Z80"
$playable.Range("E297").Value = "Z80"
$playable.Range("H297").Value = "1xAY-8910"
$playable.Range("I297").Value = "1xDAC"
$playable.Range("M297").Value = "Exciting Soccer (bootleg)"

# --- append the 1 non-working row to "GAME_NOT_WORKING FLAG" (row 5) ---
$notWorking.Range("A5").Value = 5
$notWorking.Range("B5").Value = "exctscc2"
$notWorking.Range("C5").Value = "exctsccr.c"
$notWorking.Range("D5").Value = "Z80"
$notWorking.Range("E5").Value = "Z80"
$notWorking.Range("H5").Value = "4xAY-8910"
$notWorking.Range("I5").Value = "2xDAC"
$notWorking.Range("M5").Value = "Exciting Soccer II"

# --- remove the original 4 rows (1611-1614) from "ALL"; everything below shifts up ---
$null = $all.Range("A1611:A1614").EntireRow.Delete()

# --- re-apply the autofilter over the now-smaller range ---
$all.AutoFilterMode = $false
$null = $all.Range("A1:M1747").AutoFilter()

# --- fix up the workbook-level defined names that described the old range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ALL!_FilterDatabase") {
        $n.RefersTo = "=ALL!`$A`$1:`$M`$1747"
    }
    if ($n.Name -eq "ALL!LIST") {
        $n.RefersTo = "=ALL!`$B`$1:`$M`$1747"
    }
}

# --- restore selections on each touched sheet (ALL stays the active tab) ---
$null = $playable.Range("A292:A297").Select()
$null = $notWorking.Range("A2:A5").Select()
$null = $all.Range("F1621").Select()
